$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary row 12: average of |S*|/n ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# --- New summary rows 14-17 ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# --- Formatting ---
# Build a one-off template cell fully (so only the final style is kept,
# avoiding stray intermediate cellXfs entries), then stamp its format
# onto the real target cells via PasteSpecial(Formats).

# Style used for B14:B17 -> bold, 12pt, vertically centered
$tmpl = $ws.Range("AB40")
$tmpl.Font.Bold = $true
$tmpl.Font.Size = 12
$tmpl.VerticalAlignment = -4108
$tmpl.Copy()
$ws.Range("B14:B17").PasteSpecial(-4122)
$tmpl.Clear()

# Style used for J12 -> bold, 11pt
$tmpl.Font.Bold = $true
$tmpl.Copy()
$ws.Range("J12").PasteSpecial(-4122)
$tmpl.Clear()

$ws.Range("A14:A17").EntireRow.RowHeight = 15.6

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection, matching the saved view state ---
[void]$ws.Range("A14:B17").Select()
